# Generate Report for Handback
#
# A new handback entry ("0721bd8a-ab8f-40ac-acb7-12fa72d8b191.md") is
# inserted as a new row between the existing "99ae1387..." row and the
# existing "ab4078aa..." row on all three worksheets (Overview, zh-cn,
# de-de). Each worksheet's table grows from 3 data rows to 4.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Helper: force a literal-text write (leading "'" stops the host from
# coercing date-looking / boolean-looking / numeric-looking text into a
# typed cell - it stays a plain string, matching how this workbook
# already stores even its "date" columns as shared strings).
# ---------------------------------------------------------------------
function Set-Text($range, [string]$text) {
    $range.Value = "'" + $text
}

# New GUIDs/shas used to build plausible, stable hyperlink targets for
# the newly handed-back file (mirrors the URL shape already used by the
# sibling rows - org/repo + 40-char commit sha + path).
$shaOverview = "3128f038e2da6058db33cfb8935785b0583f7f15"
$shaZhCn     = "014e4d124b627a89f4dcf53e522d7316d3124a21"
$shaDeDe     = "1b60b889f05b4bd9fd89045d0ae6602279384750"

$newFile       = "0721bd8a-ab8f-40ac-acb7-12fa72d8b191.md"
$newFileE2e    = "e2e\" + $newFile
$urlOverview   = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/$shaOverview/e2e/$newFile"
$urlZhCn       = "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/$shaZhCn/e2e/$newFile"
$urlDeDe       = "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/$shaDeDe/e2e/$newFile"

# ===== Sheet "Overview" ==================================================
$wsOverview = $wb.Worksheets.Item("Overview")

# Insert a blank row above the "ab4078aa..." row (row 3), pushing it to
# row 4. The new blank row inherits per-column formatting (hyperlink /
# date styles) from the row above, same as typical Excel insert-row.
$wsOverview.Rows.Item(3).Insert()

Set-Text $wsOverview.Range("A3") $newFile
Set-Text $wsOverview.Range("C3") ".md"
Set-Text $wsOverview.Range("E3") "Handed back: in sync with en-US"
Set-Text $wsOverview.Range("F3") "Handed back: in sync with en-US"
Set-Text $wsOverview.Range("G3") "2016-08-29 16:48:56"

# Rebuild hyperlinks from scratch - row-insert does not relocate existing
# Hyperlink objects, so the cleanest route is delete-then-readd for every
# linked cell (both the untouched rows and the two touched ones).
$wsOverview.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($wsOverview.Range("B2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c6cff3a594e900ff851748c99323ad17dba81722/e2e/99ae1387-5e58-487c-ad95-a317688eae00.md", "", "", "e2e\99ae1387-5e58-487c-ad95-a317688eae00.md") | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $urlOverview, "", "", $newFileE2e) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("B4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d7f2d0b2837c2ee8ef8a39e2e133887164867a8f/e2e/ab4078aa-f039-478f-8f90-5c6973940d4e.md", "", "", "e2e\ab4078aa-f039-478f-8f90-5c6973940d4e.md") | Out-Null

# Grow the "Overview" table (ListObject) to cover the new row.
$loOverview = $wsOverview.ListObjects.Item("Overview")
$loOverview.Resize($wsOverview.Range("A1:G4"))

# ===== Sheet "zh-cn" ======================================================
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Rows.Item(3).Insert()

Set-Text $wsZhCn.Range("A3") $newFile
Set-Text $wsZhCn.Range("B3") ".md"
Set-Text $wsZhCn.Range("C3") "Handed back: in sync with en-US"
Set-Text $wsZhCn.Range("D3") "e2e"
Set-Text $wsZhCn.Range("E3") "ht"
Set-Text $wsZhCn.Range("F3") "True"
Set-Text $wsZhCn.Range("G3") "0721bd8a-ab8f-40ac-acb7-12fa72d8b191.74b71565e6c3954bb31d95ec766c8ec4ddd55028.zh-cn.xlf"
Set-Text $wsZhCn.Range("H3") "2016-08-29 16:48:50"
Set-Text $wsZhCn.Range("I3") $newFile
Set-Text $wsZhCn.Range("J3") "0721bd8a-ab8f-40ac-acb7-12fa72d8b191.74b71565e6c3954bb31d95ec766c8ec4ddd55028.zh-cn.xlf"
Set-Text $wsZhCn.Range("K3") "2016-08-29 16:49:22"
Set-Text $wsZhCn.Range("L3") ""
Set-Text $wsZhCn.Range("M3") "ht"
Set-Text $wsZhCn.Range("N3") ""
Set-Text $wsZhCn.Range("O3") "False"
Set-Text $wsZhCn.Range("P3") ""

$wsZhCn.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c6cff3a594e900ff851748c99323ad17dba81722/e2e/99ae1387-5e58-487c-ad95-a317688eae00.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/e300a40f24e0eca84266a3521a37bcafb74ef6b5/e2e/99ae1387-5e58-487c-ad95-a317688eae00.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $urlZhCn) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $urlZhCn) | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d7f2d0b2837c2ee8ef8a39e2e133887164867a8f/e2e/ab4078aa-f039-478f-8f90-5c6973940d4e.md") | Out-Null
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-zhcn/blob/cc04364071f81128aa088d99205eaef86c5915bd/e2e/ab4078aa-f039-478f-8f90-5c6973940d4e.md") | Out-Null

$loZhCn = $wsZhCn.ListObjects.Item("zh-cn")
$loZhCn.Resize($wsZhCn.Range("A1:P4"))

# ===== Sheet "de-de" ======================================================
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Rows.Item(3).Insert()

Set-Text $wsDeDe.Range("A3") $newFile
Set-Text $wsDeDe.Range("B3") ".md"
Set-Text $wsDeDe.Range("C3") "Handed back: in sync with en-US"
Set-Text $wsDeDe.Range("D3") "e2e"
Set-Text $wsDeDe.Range("E3") "ht"
Set-Text $wsDeDe.Range("F3") "True"
Set-Text $wsDeDe.Range("G3") "0721bd8a-ab8f-40ac-acb7-12fa72d8b191.74b71565e6c3954bb31d95ec766c8ec4ddd55028.de-de.xlf"
Set-Text $wsDeDe.Range("H3") "2016-08-29 16:48:56"
Set-Text $wsDeDe.Range("I3") $newFile
Set-Text $wsDeDe.Range("J3") "0721bd8a-ab8f-40ac-acb7-12fa72d8b191.74b71565e6c3954bb31d95ec766c8ec4ddd55028.de-de.xlf"
Set-Text $wsDeDe.Range("K3") "2016-08-29 16:49:30"
Set-Text $wsDeDe.Range("L3") ""
Set-Text $wsDeDe.Range("M3") "ht"
Set-Text $wsDeDe.Range("N3") ""
Set-Text $wsDeDe.Range("O3") "False"
Set-Text $wsDeDe.Range("P3") ""

$wsDeDe.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/c6cff3a594e900ff851748c99323ad17dba81722/e2e/99ae1387-5e58-487c-ad95-a317688eae00.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/d9c0a3b6d77e4fcbb2fdbb2437f06cad9bf09fbf/e2e/99ae1387-5e58-487c-ad95-a317688eae00.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $urlDeDe) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $urlDeDe) | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A4"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d7f2d0b2837c2ee8ef8a39e2e133887164867a8f/e2e/ab4078aa-f039-478f-8f90-5c6973940d4e.md") | Out-Null
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I4"), "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/525516b7b8991a834e5fabdc5e5b03bd8cc9a238/e2e/ab4078aa-f039-478f-8f90-5c6973940d4e.md") | Out-Null

$loDeDe = $wsDeDe.ListObjects.Item("de-de")
$loDeDe.Resize($wsDeDe.Range("A1:P4"))

Write-Host "Done: inserted 0721bd8a row on Overview/zh-cn/de-de sheets."
